$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row to append: row 84, date serial 45884 (2025-08-15) and value -0.4792091214565772
$newRow = 84

# Set the date value in column A, copying the style (incl. date number format) from the row above
$ws.Cells.Item($newRow - 1, 1).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item($newRow, 1).Value = 45884
$ws.Cells.Item($newRow, 2).Value = -0.4792091214565772
